$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Positive / Negative / Neutral tweet-sentiment counts for each movie row
# (columns G, H, I of Table2), rows 2-16.
$ws.Cells.Item(2, 7).Value = 60232
$ws.Cells.Item(2, 8).Value = 29703
$ws.Cells.Item(2, 9).Value = 56955

$ws.Cells.Item(3, 7).Value = 58112
$ws.Cells.Item(3, 8).Value = 83316
$ws.Cells.Item(3, 9).Value = 26402

$ws.Cells.Item(4, 7).Value = 67037
$ws.Cells.Item(4, 8).Value = 16596
$ws.Cells.Item(4, 9).Value = 122824

$ws.Cells.Item(5, 7).Value = 48126
$ws.Cells.Item(5, 8).Value = 10434
$ws.Cells.Item(5, 9).Value = 58383

$ws.Cells.Item(6, 7).Value = 38070
$ws.Cells.Item(6, 8).Value = 9163
$ws.Cells.Item(6, 9).Value = 73528

$ws.Cells.Item(7, 7).Value = 59583
$ws.Cells.Item(7, 8).Value = 17456
$ws.Cells.Item(7, 9).Value = 109801

$ws.Cells.Item(8, 7).Value = 44821
$ws.Cells.Item(8, 8).Value = 7417
$ws.Cells.Item(8, 9).Value = 60184

$ws.Cells.Item(9, 7).Value = 23968
$ws.Cells.Item(9, 8).Value = 4762
$ws.Cells.Item(9, 9).Value = 42884

$ws.Cells.Item(10, 7).Value = 30982
$ws.Cells.Item(10, 8).Value = 8036
$ws.Cells.Item(10, 9).Value = 85365

$ws.Cells.Item(11, 7).Value = 22707
$ws.Cells.Item(11, 8).Value = 3498
$ws.Cells.Item(11, 9).Value = 36183

$ws.Cells.Item(12, 7).Value = 12572
$ws.Cells.Item(12, 8).Value = 2061
$ws.Cells.Item(12, 9).Value = 18354

$ws.Cells.Item(13, 7).Value = 8751
$ws.Cells.Item(13, 8).Value = 3268
$ws.Cells.Item(13, 9).Value = 10782

$ws.Cells.Item(14, 7).Value = 4600
$ws.Cells.Item(14, 8).Value = 1373
$ws.Cells.Item(14, 9).Value = 4268

$ws.Cells.Item(15, 7).Value = 2385
$ws.Cells.Item(15, 8).Value = 690
$ws.Cells.Item(15, 9).Value = 4560

$ws.Cells.Item(16, 7).Value = 2275
$ws.Cells.Item(16, 8).Value = 484
$ws.Cells.Item(16, 9).Value = 3454

# Mirror the reviewer's final on-screen state: scrolled right so column H is
# the first visible (unfrozen) column, with I4 as the active selection.
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("I4").Select()
